# Apply edit: insert 3 new rows (762,763,764) into the Plátano sheet,
# pushing the previous rows 762..860 down to 765..863.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows starting at row 762 (shifts existing rows down)
$ws.Rows("762:764").Insert()

# Column layout:
# A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Tipo,
# G Producto ID, H Producto, I Categoría ID, J Categoría,
# K Variedad, L Calidad, M Volumen, N Precio mínimo, O Precio máximo,
# P Precio promedio ponderado, Q Unidad de comercialización, R Origen,
# S Precio $/Kg, T Kg / unidad

# Row 762
$ws.Cells.Item(762, 1).Value = 10
$ws.Cells.Item(762, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(762, 3).Value = "La Araucanía"
$ws.Cells.Item(762, 4).Value = 44918
$ws.Cells.Item(762, 5).Value = 9
$ws.Cells.Item(762, 6).Value = "Fruta"
$ws.Cells.Item(762, 7).Value = 100108
$ws.Cells.Item(762, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(762, 9).Value = 100108006
$ws.Cells.Item(762, 10).Value = "Plátano"
$ws.Cells.Item(762, 11).Value = "Barraganete"
$ws.Cells.Item(762, 12).Value = "Maduro"
$ws.Cells.Item(762, 13).Value = 50
$ws.Cells.Item(762, 14).Value = 50000
$ws.Cells.Item(762, 15).Value = 50000
$ws.Cells.Item(762, 16).Value = 50000
$ws.Cells.Item(762, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(762, 18).Value = "Ecuador"
$ws.Cells.Item(762, 19).Value = 2500
$ws.Cells.Item(762, 20).Value = 20

# Row 763
$ws.Cells.Item(763, 1).Value = 10
$ws.Cells.Item(763, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(763, 3).Value = "La Araucanía"
$ws.Cells.Item(763, 4).Value = 44918
$ws.Cells.Item(763, 5).Value = 9
$ws.Cells.Item(763, 6).Value = "Fruta"
$ws.Cells.Item(763, 7).Value = 100108
$ws.Cells.Item(763, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(763, 9).Value = 100108006
$ws.Cells.Item(763, 10).Value = "Plátano"
$ws.Cells.Item(763, 11).Value = "Barraganete"
$ws.Cells.Item(763, 12).Value = "Verde"
$ws.Cells.Item(763, 13).Value = 50
$ws.Cells.Item(763, 14).Value = 48000
$ws.Cells.Item(763, 15).Value = 48000
$ws.Cells.Item(763, 16).Value = 48000
$ws.Cells.Item(763, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(763, 18).Value = "Ecuador"
$ws.Cells.Item(763, 19).Value = 2400
$ws.Cells.Item(763, 20).Value = 20

# Row 764
$ws.Cells.Item(764, 1).Value = 10
$ws.Cells.Item(764, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(764, 3).Value = "La Araucanía"
$ws.Cells.Item(764, 4).Value = 44918
$ws.Cells.Item(764, 5).Value = 9
$ws.Cells.Item(764, 6).Value = "Fruta"
$ws.Cells.Item(764, 7).Value = 100108
$ws.Cells.Item(764, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(764, 9).Value = 100108006
$ws.Cells.Item(764, 10).Value = "Plátano"
$ws.Cells.Item(764, 11).Value = "Sin especificar"
$ws.Cells.Item(764, 12).Value = "Pintón"
$ws.Cells.Item(764, 13).Value = 550
$ws.Cells.Item(764, 14).Value = 22000
$ws.Cells.Item(764, 15).Value = 23000
$ws.Cells.Item(764, 16).Value = 22545
$ws.Cells.Item(764, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(764, 18).Value = "Ecuador"
$ws.Cells.Item(764, 19).Value = 1127
$ws.Cells.Item(764, 20).Value = 20

# Make sure the Fecha (date) column keeps the existing date number format
$ws.Range("D762:D764").NumberFormat = $ws.Range("D765").NumberFormat
